# 6.1.1.xlsx update: add a new "2020" data column (M) to the indicator table,
# matching the formatting of the existing 2019 column (L), and move the
# active selection to G15 (with the view scrolled so column B is leftmost).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlPasteFormats = -4122

function Copy-ColumnLFormat {
    param([string]$row)
    $ws.Range("L$row").Copy() | Out-Null
    $ws.Range("M$row").PasteSpecial($xlPasteFormats) | Out-Null
    $excel.CutCopyMode = $false
}

# Row 2 - thin separator row under the title, same style as the rest of the row.
Copy-ColumnLFormat "2"

# Row 3 - header row: new year column, 2020.
Copy-ColumnLFormat "3"
$ws.Range("M3").Value = 2020

# Row 4 - national ("Кыргызская Республика") total row: new value is shown in
# bold to highlight it (a brand-new cell style, same number format/border as
# the rest of the row but with a bold font).
Copy-ColumnLFormat "4"
$ws.Range("M4").Value = 94.1
$ws.Range("M4").Font.Bold = $true

# Remaining data rows - plain new values, same formatting as column L.
Copy-ColumnLFormat "5"
$ws.Range("M5").Value = 99.6

Copy-ColumnLFormat "6"
$ws.Range("M6").Value = 91

Copy-ColumnLFormat "7"
$ws.Range("M7").Value = 86.886172668979881

Copy-ColumnLFormat "8"
$ws.Range("M8").Value = 86.955790296225956

Copy-ColumnLFormat "9"
$ws.Range("M9").Value = 96.29195112324031

Copy-ColumnLFormat "10"
$ws.Range("M10").Value = 97.849780305474511

Copy-ColumnLFormat "11"
$ws.Range("M11").Value = 90.676703333930902

Copy-ColumnLFormat "12"
$ws.Range("M12").Value = 99.675929342188979

Copy-ColumnLFormat "13"
$ws.Range("M13").Value = 100

Copy-ColumnLFormat "14"
$ws.Range("M14").Value = 100

Copy-ColumnLFormat "15"
$ws.Range("M15").Value = 100

# Update the view: scroll so column B is the leftmost visible column, and
# leave the active selection on G15.
$win = $excel.ActiveWindow
$ws.Range("G15").Select() | Out-Null
$win.ScrollColumn = 2
